$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the typo/formatting bug: the note text was missing a trailing '*'
# (gdata:::read.xlsx mis-read the quoted text). Append the missing
# asterisk to restore '***soon***'.
$ws.Range("C2").Value = 'An **EXAMPLE**; feel "free" to *delete* it ***soon***'

# Move the active selection to the corrected cell.
[void]$ws.Range("C2").Select()
